# Generate Report for Handoff
# The 6e4d3395-ef20-4ad4-9f46-cdb09eaec372 entry has been handed off again,
# so its "Latest Handoff" timestamps are refreshed on every localized sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column for the 6e4d3395... row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-41-13 10:41:48"

# --- zh-cn sheet: "Latest Handoff Datetime" column for the 6e4d3395... row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-13 10:41:45"

# --- de-de sheet: "Latest Handoff Datetime" column for the 6e4d3395... row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-13 10:41:48"
